$d = $word.ActiveDocument

# --- 1. "hustle" -> "hassle" -----------------------------------------
$d.Content.Find.Execute("hustle", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "hassle", 2) | Out-Null

# --- 2. "MOULES:" -> "MODULES:" ---------------------------------------
$d.Content.Find.Execute("MOULES:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MODULES:", 2) | Out-Null

# --- 3. "RELATIONSHIPS/COMMUNICATION:" -> "RELATIONSHIPS / COMMUNICATION:"
$d.Content.Find.Execute("RELATIONSHIPS/COMMUNICATION:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "RELATIONSHIPS / COMMUNICATION:", 2) | Out-Null

# --- 4. FUNCTIONS / MODULES list: insert a new "Signup" entry right
#        after "Login", shifting the remaining entries down by one and
#        dropping the trailing "/FAQ" from the last item. ---------------
# Paragraph 22 originally holds two runs ("Online support" + "/FAQ"); a
# plain Range.Text assignment only touches the first run, so delete the
# whole paragraph content first, then insert the replacement text.
$rVolunteering = $d.Paragraphs(22).Range
$rVolunteering.End = $rVolunteering.End - 1
$rVolunteering.Delete()
$rVolunteering = $d.Paragraphs(22).Range
$rVolunteering.End = $rVolunteering.End - 1
$rVolunteering.InsertAfter("Volunteering")

$d.Paragraphs(21).Range.Text = "Payment"
$d.Paragraphs(20).Range.Text = "Participation"
$d.Paragraphs(19).Range.Text = "Search"
$d.Paragraphs(18).Range.Text = "Notification Service"
$d.Paragraphs(17).Range.Text = "Event"
$d.Paragraphs(16).Range.Text = "Signup"

# --- 5. Move the "_GoBack" bookmark from the "event management"
#        paragraph onto the end of the new "Signup" entry. --------------
$rSignup = $d.Paragraphs(16).Range
$rSignup.End = $rSignup.End - 1
$d.Bookmarks.Add("_GoBack", $rSignup) | Out-Null
